# Add a new test user "a.madison@microfocus.com" to the Global sheet's user list,
# mirroring the formatting of the row above it, and restore the originally
# active sheet/tab afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Copy formatting (style) from the last existing data row (A5) into the new row (A6)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Set the new user's email address
$ws.Range("A6").Value = "a.madison@microfocus.com"

# Update the selection on the Global sheet to the newly added cell
$ws.Range("A6").Select()

# Restore the workbook's originally active sheet/tab (ChangePW)
$wb.Worksheets.Item("ChangePW").Activate()
